$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the User Name / Password / Confirm Password for Charlie's row (row 4)
$ws.Range("C4").Value = "testcypress"
$ws.Range("D4").Value = "ValidPass123"
$ws.Range("E4").Value = "DifferentP123"

# Remove the last data row (row 5 - John Doe) entirely
$ws.Rows(5).Delete()

# Move the selection to where the author last left the cursor
$ws.Range("D8").Select() | Out-Null
